# "Common: Added vape edit" -- append translation rows for the new
# lab.vape.edit.* / lab.vape.update.* / lab.vape.index.* keys to the
# "Import" sheet (sheet1), mirroring the existing lab.mixture.* /
# lab.build.* rows already present.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# New translation rows to append right after the last existing row (519).
# Each entry is: key, translation (cs)
$rows = @(
    @("lab.vape.edit.title", "Editace vapu"),
    @("lab.vape.edit.subtitle", "Každý se někdy překlepne, zde je možné upravit vape."),
    @("lab.vape.link.button", "Zpět"),
    @("lab.vape.update.submit", "Aktualizovat"),
    @("lab.vape.update.success", "Vape byl aktualizován."),
    @("lab.vape.index.title", "Náhled vapu"),
    @("lab.vape.button.edit", "Editovat"),
    @("lab.vape.index.preview.subtitle", "Správa vybraného vapu"),
    @("lab.vape.index.preview.title", "Náhled vapu")
)

$lastRow = 519
$startRow = $lastRow + 1

for ($i = 0; $i -lt $rows.Count; $i++) {
    $targetRow = $startRow + $i
    # Duplicate the last row (same "cs" language + formatting/style) and
    # insert it right below the previously inserted row.
    $ws.Rows($lastRow).Copy()
    $ws.Rows($targetRow).Insert()

    $ws.Range("A$targetRow").Value = "cs"
    $ws.Range("B$targetRow").Value = $rows[$i][0]
    $ws.Range("C$targetRow").Value = $rows[$i][1]
}

# Match the saved selection/scroll position from the authored workbook.
$ws.Range("B524").Select()
